# Remove the "Number of Promotions" row (row 4) and the
# "Number of Dependents" row (originally row 6, now row 5 after the
# first deletion), leaving "Distance from Home" as the last row (row 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(4).Delete()
$ws.Rows.Item(5).Delete()
